$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.567.83'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.938.22'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").Value = '2.939.43'
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.439'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.32%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '3.479.12'
$ws.Range("E13").Value = '  +2.44%  '
$ws.Range("E14").Value = '  -1.27%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.83%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '76.496.72'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '2.951.53'
$ws.Range("E18").Value = '  +2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.07%  '
$ws.Range("E20").Value = '  -3.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '373.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = '3.094.82'
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000106'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.25%  '
$ws.Range("E32").Value = '  -3.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '496.72'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.44%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.74%  '
$ws.Range("E37").Value = '  +20.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  +13.01%  '
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -6.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '179.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").Value = '  -4.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.588'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.89'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("E51").Value = '  -5.96%  '
